# Update "想去人数" (people interested) counts in the workbook.
# Sheet "展览" (sheet1)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 649
$ws1.Range("F3").Value = 737
$ws1.Range("F4").Value = 954
$ws1.Range("F5").Value = 735
$ws1.Range("F6").Value = 845
$ws1.Range("F7").Value = 410
$ws1.Range("F16").Value = 597
$ws1.Range("F17").Value = 3
$ws1.Range("F18").Value = 364
$ws1.Range("F23").Value = 589
$ws1.Range("F25").Value = 809

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 330
$ws2.Range("F6").Value = 23

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 649
$ws4.Range("F6").Value = 330
$ws4.Range("F7").Value = 737
$ws4.Range("F8").Value = 954
$ws4.Range("F9").Value = 735
$ws4.Range("F10").Value = 845
$ws4.Range("F11").Value = 410
$ws4.Range("F17").Value = 23
$ws4.Range("F19").Value = 518
$ws4.Range("F23").Value = 597
$ws4.Range("F25").Value = 3
$ws4.Range("F26").Value = 364
$ws4.Range("F37").Value = 589
$ws4.Range("F39").Value = 809
